# Update Marcus Stoinis / Delhi Capitals innings stats (runs, balls, fours, sixes)
# Values are written as text (matching the original <c t="str"> storage) by
# forcing a text number format before assigning the string values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ C = "53"; D = "26"; E = "6"; F = "2" }
    3  = @{ C = "1";  D = "3";  E = "0"; F = "0" }
    4  = @{ C = "0";  D = "1";  E = "0"; F = "0" }
    7  = @{ C = "2";  D = "3";  E = "0"; F = "0" }
    8  = @{ C = "9";  D = "10"; E = "0"; F = "0" }
    9  = @{ C = "10"; D = "5";  E = "0"; F = "1" }
    11 = @{ C = "5";  D = "3";  E = "1"; F = "0" }
    12 = @{ C = "11"; D = "9";  E = "1"; F = "0" }
    13 = @{ C = "53"; D = "21"; E = "7"; F = "3" }
    14 = @{ C = "24"; D = "14"; E = "1"; F = "2" }
    15 = @{ C = "18"; D = "19"; E = "1"; F = "0" }
    16 = @{ C = "5";  D = "6";  E = "1"; F = "0" }
    17 = @{ C = "6";  D = "6";  E = "0"; F = "0" }
}

foreach ($row in $updates.Keys) {
    $rowValues = $updates[$row]
    foreach ($col in $rowValues.Keys) {
        $cell = $ws.Range("$col$row")
        $cell.NumberFormat = "@"
        $cell.Value = $rowValues[$col]
    }
}
